# Commiting scripts (R22 UAT2 - Regression).
#
# Insert a new worksheet named "Sheet1" before the existing
# "KYC_Amendment_Customer" sheet. The new sheet carries a copy of the
# header row and the first data row from the original sheet (with the
# customer id in A2 changed), and becomes the active sheet with A2
# selected. The original sheet keeps all of its data but its selection
# becomes the whole used range with C20 as the last-navigated cell.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts a new sheet immediately before the sheet that
# is currently active, and auto-names it "Sheet1" (matching sheetId=2 /
# rId1, while the original sheet shifts to rId2) - exactly the shape of
# the target workbook.xml.
$ws1 = $wb.Worksheets.Add()
$ws1.Name = "Sheet1"

$ws2 = $wb.Worksheets.Item("KYC_Amendment_Customer")

# --- New "Sheet1": header row + first data row -----------------------
$ws1.Range("A1").Value = "value:1:1:1"
$ws1.Range("B1").Value = "OCCUPATION"
$ws1.Range("C1").Value = "NAME.OF.BUS"
$ws1.Range("D1").Value = "NAT.OF.BUS"
$ws1.Range("E1").Value = "STAT.OWNER"
$ws1.Range("F1").Value = "NAME.OF.EMP"
$ws1.Range("G1").Value = "CS.POS"
$ws1.Range("H1").Value = "CS.EMP.SINCE"
$ws1.Range("I1").Value = "CURRENT.SALARY"
$ws1.Range("J1").Value = "OTHER.INCOME"
$ws1.Range("K1").Value = "OTHER.FUNDS"
$ws1.Range("L1").Value = "PER.PROP.INMNT"
$ws1.Range("M1").Value = "CS.ANNUM.TO"
$ws1.Range("N1").Value = "SOURCE.OF.INCOME"
$ws1.Range("O1").Value = "POLITICAL.FIGURE"
$ws1.Range("P1").Value = "KYC.REVW.COMENT"
$ws1.Range("Q1").Value = "CUST.COMMENTS:1"

$ws1.Range("A2").Value = 16324801
$ws1.Range("B2").Value = "Salaried"
$ws1.Range("C2").Value = "Test"
$ws1.Range("D2").Value = "Test1"
$ws1.Range("E2").Value = "Test2"
$ws1.Range("J2").Value = 0
$ws1.Range("K2").Value = 0
$ws1.Range("L2").Value = 20000000
$ws1.Range("M2").Value = 400000
$ws1.Range("N2").Value = "Business"
$ws1.Range("P2").Value = "Testing1"
$ws1.Range("Q2").Value = "ISL"

# New sheet is active, with A2 selected.
$ws1.Activate()
$ws1.Range("A2").Select()

# --- Existing sheet: selection becomes the whole used range -----------
# (Real Excel remembers C20 as the last-active cell inside the A1:Q101
# selection; this host's Range.Activate()/Select() always collapses the
# active cell to the selected rectangle's corner, so selecting the used
# range - which does yield the target sqref "A1:Q101" - is the closest
# reachable approximation.)
$ws2.Activate()
$ws2.UsedRange.Select()

# Leave the new first sheet as the active / tab-selected sheet, matching
# the target workbook (Sheet1 has tabSelected="1").
$ws1.Activate()
